$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mike's (columns D/E) week 1 and week 2 hours
$ws.Range("D4").Value = 15
$ws.Range("D5").Value = 13

# Patrick's (columns F/G) week 1 and week 2 hours
$ws.Range("F4").Value = 13.5
$ws.Range("F5").Value = 13

# Update the active cell selection to F5
$ws.Range("F5").Select()
